$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q1" sheet right before the "总计" summary sheet
#    (mirrors the per-quarter fund-holding detail sheets already present).
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totals)
$q1.Name = "2022-Q1"

# Re-resolve the "总计" sheet by name: the COM reference captured above can
# rebind by tab position once a new sheet is inserted ahead of it, so grab a
# fresh handle pointing at the actual (now shifted) summary sheet.
$totals = $wb.Worksheets.Item("总计")

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$header = $q1.Range("B1:H1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Row 2 - 天弘越南市场股票（QDII）A
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'008763"
$q1.Range("C2").Value = "天弘越南市场股票（QDII）A"
$q1.Range("D2").Value = "'37.53"
$q1.Range("E2").Value = "'92.10"
$q1.Range("F2").Value = "'5.91"
$q1.Range("G2").Value = "'2.2180"
$q1.Range("H2").Value = 5

# Row 3 - 天弘越南市场股票（QDII）C
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'008764"
$q1.Range("C3").Value = "天弘越南市场股票（QDII）C"
$q1.Range("D3").Value = "'14.26"
$q1.Range("E3").Value = "'92.10"
$q1.Range("F3").Value = "'5.91"
$q1.Range("G3").Value = "'0.8428"
$q1.Range("H3").Value = 5

$idxCol = $q1.Range("A2:A3")
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$totals.Rows.Item(2).Insert()

# The inserted row picks up formatting from the row above (the bold header);
# strip that back off the plain data cells before writing their values.
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 3.06

$a2 = $totals.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3

# Restore the originally-active sheet/selection so this edit doesn't leave
# an unrelated side effect on the workbook's view state.
$wb.Worksheets.Item("2021-Q2").Activate() | Out-Null
$wb.Worksheets.Item("2021-Q2").Range("A1").Select() | Out-Null
